$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells with refreshed media_hero URLs (I117, I124) ---
# --- and fill in previously-empty media_hero/media_alt cells (I128/J128, I130/J130) ---
$ws.Range("I117").Value = "https://drive.google.com/file/d/16lshwnpjHOoXJAHhSGtCTLvguDmVU8Zg/view?usp=drive_link"
$ws.Range("I124").Value = "https://drive.google.com/file/d/1BxWu91KBh1TCM6eJsotpHkoyV-jOYtSV/view?usp=drive_link"
$ws.Range("I128").Value = "https://drive.google.com/file/d/1dAhSBpGPjazW0Ok_tlYhaCd6KjnRCxZa/view?usp=drive_link"
$ws.Range("J128").Value = "Self-Destructive, Soul Feedback"
$ws.Range("I130").Value = "https://drive.google.com/file/d/1irw5VYRZ4kkm5XhLER1FB6CVCxjmaSX8/view?usp=drive_link"
$ws.Range("J130").Value = "Armored Heart, Soul Feedback"

# --- Append 5 new rows (140-144) generated from the Google Sheet source ---
# Row 140: Concealed
$ws.Range("A140").Value = "songs"
$ws.Range("B140").Value = "Concealed"
$ws.Range("C140").Value = "concealed-plastic-extensions-redemptions"
$ws.Range("D140").Value = "/songs/concealed-plastic-extensions-redemptions/"
$ws.Range("E140").Value = "MusicComposition"
$ws.Range("H140").Value = "concealing, plastic, extensions, redemptions, con artist, framming"
$ws.Range("I140").Value = "https://drive.google.com/file/d/18_tupZGmpdqRNwjyQ8Wcof9SVKCwgxuC/view?usp=drive_link"
$ws.Range("J140").Value = "Concealed, Soul Feedback"
$ws.Range("K140").Value = "Body Dismorphic Disorder"
$ws.Range("L140").Value = "Collateral"
$ws.Range("N140").Value = "CONCEALED`nJust passing by`nand Framing`n(+) Filter-ing.`nConcealer,`nwith plastic,`nredemptions,`nextensions.`nSealed tight,`nCon artist.`nIn decades when that make up`nfades away`nGhost-portraits:`nsomeone`nwho wasn't never really even there.`nConcealed`nBuying,`nBorrowing,`nRenting."

# Row 141: Dance
$ws.Range("A141").Value = "songs"
$ws.Range("B141").Value = "Dance"
$ws.Range("C141").Value = "dance-electrocardiogram-beauty-eclipse-midday"
$ws.Range("D141").Value = "/songs/dance-electrocardiogram-beauty-eclipse-midday/"
$ws.Range("E141").Value = "MusicComposition"
$ws.Range("H141").Value = "dance, electrocardiogram, beauty, midday, eclipse"
$ws.Range("I141").Value = "https://drive.google.com/file/d/1N4CgRx5ooNFhqKhchcg4VW8yN4fnEQre/view?usp=drive_link"
$ws.Range("J141").Value = "Dance, Soul Feedback"
$ws.Range("L141").Value = "Sunk"
$ws.Range("N141").Value = "DANCE`n[Part I Bossa Nova]`nDance,`nthe sea`nwill survive`nour rush.`nDance,`nWe only remember peaks`namong the waves of electrocardiograms.`nDance,`nthat beauty is a eclipse`nat midday`nAnd it'll pass through us.`n[Part II Heavy Instrumental Descarga]`n[Part III With Feeling slow Son Cubano]`nYou make all this sea`nYou and this city`nIf you're not here`nthere's nothing at all...`nIf you're not here`nThe evil in me..."

# Row 142: Dejavu
$ws.Range("A142").Value = "songs"
$ws.Range("B142").Value = "Dejavu"
$ws.Range("C142").Value = "dejavu-predestination-repetition-compulsively-commitment"
$ws.Range("D142").Value = "/songs/dejavu-predestination-repetition-compulsively-commitment/"
$ws.Range("E142").Value = "MusicComposition"
$ws.Range("H142").Value = "dejavu, predestination, repetition, compulsively, commitment"
$ws.Range("I142").Value = "https://drive.google.com/file/d/1hgo8yjtBI8ol-rbOxEEkokqb9Nr5xThz/view?usp=drive_link"
$ws.Range("J142").Value = "Dejavu, Soul Feedback"
$ws.Range("L142").Value = "Animal"
$ws.Range("N142").Value = "DEJAVU`nYou talk about love`nAnd predestination`nBut every beginning`nIs just repetition`nAnd that's no love`nThat's dejavu.`nCompulsively avoiding commitment`nThat's to return.`nAnd every time is less intense,`nEvery time you need more to feel the same,`nThe same or less.`nSo don't get scared`nIf after moaning bodies comes destruction."

# Row 143: Dissonance
$ws.Range("A143").Value = "songs"
$ws.Range("B143").Value = "Dissonance"
$ws.Range("C143").Value = "dissonance-doublesidedknife-memory-redeem-rage"
$ws.Range("D143").Value = "/songs/dissonance-doublesidedknife-memory-redeem-rage/"
$ws.Range("E143").Value = "MusicComposition"
$ws.Range("H143").Value = "dissonance, doublesidedknife, memory, redeem, rage"
$ws.Range("I143").Value = "https://drive.google.com/file/d/1F9KRRjwUM7Ukb6YrK025lmMF2CHabllw/view?usp=drive_link"
$ws.Range("J143").Value = "Dissonace, Soul Feedback"
$ws.Range("L143").Value = "Counterfeit"
$ws.Range("N143").Value = "DISSONANCE`nShe searches for a double sided knife,`nTo build a memory.`nHe said bring your lover to our house,`nAnd so we all get to feel the same.`nA painful promise,`nKeeps the thread of the story,`nTo redeem the animal-rage:`nDissonance.`nThe idea of giving birth`nDon't compensate the hurting`nbrought everyday by a dying love.`n`"I feel like I wish you had someone you liked`""

# Row 144: Granade
$ws.Range("A144").Value = "sculpture"
$ws.Range("B144").Value = "Granade"
$ws.Range("C144").Value = "granade-heart-fire-relationship-timer-flammable"
$ws.Range("D144").Value = "/sculpture/granade-heart-fire-relationship-timer-flammable/"
$ws.Range("E144").Value = "Sculpture"
$ws.Range("H144").Value = "granade, heart, fire, relationship"
$ws.Range("I144").Value = "https://drive.google.com/file/d/1a08h8JyNOQNB2kFHhZZxguyUC19fKWyu/view?usp=drive_link"
$ws.Range("J144").Value = "Granade, Soul Feedback. Zodiac of Soullessness, Glass sculpture"
$ws.Range("K144").Value = "Granade, Heart"
$ws.Range("N144").Value = "`"Foreseeing fire in your relationships`""

